$wb = $excel.ActiveWorkbook

# --- Sheet "Metadata" updates ---
$meta = $wb.Worksheets.Item("Metadata")

# Version 5.0.0 -> 6.0.0
$meta.Range("B3").Value = "6.0.0"

# Date update
$meta.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher value was blank -> "Alvearie Team"
$meta.Range("B9").Value = "Alvearie Team"

# Row 10 was "Contact" / "No display for ContactDetail" -> "Jurisdiction" / "United States of America"
$meta.Range("A10").Value = "Jurisdiction"
$meta.Range("B10").Value = "United States of America"

# Row 11 was a duplicate "Contact" / "No display for ContactDetail" row -> remove entirely
$meta.Rows("11:11").Delete()

# --- Sheet "Elements" updates ---
$elements = $wb.Worksheets.Item("Elements")

# Root extension row: Short/Definition change from generic placeholders to specific text
$elements.Range("K2").Value = "Reimburse Type"
$elements.Range("L2").Value = "Method of payment code for the claim"
